$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.612.24"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "3.308.51"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'518.08"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'170.74"
$ws.Range("E6").Value = "  -6.63%  "

$ws.Range("D7").Value = "'0.587"
$ws.Range("E7").Value = "  -2.79%  "

$ws.Range("D8").Value = "3.291.53"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.599"
$ws.Range("E10").Value = "  -4.00%  "

$ws.Range("D11").Value = "'52.52"
$ws.Range("E11").Value = "  -12.29%  "

$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -2.55%  "

$ws.Range("D14").Value = "'8.91"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").Value = "3.866.70"
$ws.Range("E15").Value = "  +1.59%  "

$ws.Range("D16").Value = "3.325.69"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").Value = "'0.116"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.478.20"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.34"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").Value = "'11.07"
$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("D21").Value = "'0.948"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "'371.99"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  +7.05%  "

$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("D25").Value = "'80.98"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").Value = "'3.63"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("E27").Value = "  +1.44%  "

$ws.Range("D28").Value = "'2.65"
$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").Value = "'11.12"
$ws.Range("E29").Value = "  -4.08%  "

$ws.Range("D30").Value = "'8.08"
$ws.Range("E30").Value = "  -4.28%  "

$ws.Range("D31").Value = "'28.52"
$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("D32").Value = "'621.18"
$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").Value = "'6.32"
$ws.Range("E33").Value = "  -9.10%  "

$ws.Range("D34").Value = "'11.09"
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("E35").Value = "  -2.52%  "

$ws.Range("D36").Value = "'57.57"
$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "'35.49"
$ws.Range("E38").Value = "  -4.27%  "

$ws.Range("D39").Value = "'0.373"
$ws.Range("E39").Value = "  -7.83%  "

$ws.Range("D40").Value = "0.0₃0718"
$ws.Range("E40").Value = "  +6.80%  "

$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").Value = "'2.60"
$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").Value = "2.886.39"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("D44").Value = "'0.122"
$ws.Range("E44").Value = "  -3.87%  "

$ws.Range("D45").Value = "'2.99"
$ws.Range("E45").Value = "  +1.85%  "

$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("D47").Value = "'0.0391"
$ws.Range("E47").Value = "  -2.49%  "

$ws.Range("D48").Value = "'2.56"
$ws.Range("E48").Value = "  -6.65%  "

$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'135.44"
$ws.Range("E50").Value = "  +2.39%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.123"
$ws.Range("E51").Value = "  -2.18%  "
